$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.167.99'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '3.524.32'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.47'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.13'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').Value = '3.524.73'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.14'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.378'
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('D13').Value = '4.125.40'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.48'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.118'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '3.525.66'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '64.205.46'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.80'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.92'
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '382.79'
$ws.Range('E22').Value = '  -2.31%  '
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('D24').Value = '3.666.51'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('E28').Value = '  +3.52%  '
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.47'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.44'
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').Value = '3.538.33'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.61'
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.146'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.42'
$ws.Range('E38').Value = '  +3.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.95'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '160.49'
$ws.Range('E41').Value = '  -5.18%  '
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.69'
$ws.Range('E43').Value = '  +2.37%  '
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.41'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('D50').Value = '2.476.34'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('E51').Value = '  -0.98%  '
